$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = "34.489.82"
$ws.Cells.Item(2, 5).Value = "  +0.44%  "
$ws.Cells.Item(3, 4).Value = "1.810.73"
$ws.Cells.Item(3, 5).Value = "  +0.44%  "
$ws.Cells.Item(4, 5).Value = "  -0.11%  "
$ws.Cells.Item(5, 4).NumberFormat = "@"
$ws.Cells.Item(5, 4).Value = "225.48"
$ws.Cells.Item(5, 5).Value = "  -1.00%  "
$ws.Cells.Item(6, 4).NumberFormat = "@"
$ws.Cells.Item(6, 4).Value = "0.593"
$ws.Cells.Item(6, 5).Value = "  +2.66%  "
$ws.Cells.Item(8, 4).NumberFormat = "@"
$ws.Cells.Item(8, 4).Value = "38.24"
$ws.Cells.Item(8, 5).Value = "  +5.74%  "
$ws.Cells.Item(9, 5).Value = "  -4.27%  "
$ws.Cells.Item(10, 5).Value = "  -2.67%  "
$ws.Cells.Item(11, 4).NumberFormat = "@"
$ws.Cells.Item(11, 4).Value = "0.0974"
$ws.Cells.Item(11, 5).Value = "  +0.93%  "
$ws.Cells.Item(12, 4).Value = "2.072.24"
$ws.Cells.Item(12, 5).Value = "  +0.48%  "
$ws.Cells.Item(13, 5).Value = "  -3.10%  "
$ws.Cells.Item(14, 4).Value = "1.801.39"
$ws.Cells.Item(14, 5).Value = "  -0.29%  "
$ws.Cells.Item(15, 4).NumberFormat = "@"
$ws.Cells.Item(15, 4).Value = "0.632"
$ws.Cells.Item(15, 5).Value = "  -1.99%  "
$ws.Cells.Item(16, 4).Value = "34.455.25"
$ws.Cells.Item(17, 5).Value = "  -1.69%  "
$ws.Cells.Item(18, 4).NumberFormat = "@"
$ws.Cells.Item(18, 4).Value = "68.30"
$ws.Cells.Item(18, 5).Value = "  -1.15%  "
$ws.Cells.Item(19, 4).NumberFormat = "@"
$ws.Cells.Item(19, 4).Value = "243.08"
$ws.Cells.Item(19, 5).Value = "  -1.07%  "
$ws.Cells.Item(20, 4).Value = "0.0₃0773"
$ws.Cells.Item(20, 5).Value = "  -2.80%  "
$ws.Cells.Item(21, 4).NumberFormat = "@"
$ws.Cells.Item(21, 4).Value = "11.21"
$ws.Cells.Item(21, 5).Value = "  -2.65%  "
$ws.Cells.Item(22, 5).Value = "  -0.07%  "
$ws.Cells.Item(23, 5).Value = "  -1.40%  "
$ws.Cells.Item(24, 5).Value = "  +3.41%  "
$ws.Cells.Item(25, 4).NumberFormat = "@"
$ws.Cells.Item(25, 4).Value = "170.20"
$ws.Cells.Item(25, 5).Value = "  -1.21%  "
$ws.Cells.Item(26, 4).NumberFormat = "@"
$ws.Cells.Item(26, 4).Value = "7.80"
$ws.Cells.Item(26, 5).Value = "  -1.71%  "
$ws.Cells.Item(27, 4).NumberFormat = "@"
$ws.Cells.Item(27, 4).Value = "17.60"
$ws.Cells.Item(27, 5).Value = "  +4.39%  "
$ws.Cells.Item(28, 5).Value = "  +1.59%  "
$ws.Cells.Item(29, 5).Value = "  -0.08%  "
$ws.Cells.Item(30, 5).Value = "  -1.42%  "
$ws.Cells.Item(31, 5).Value = "  -1.12%  "
$ws.Cells.Item(32, 4).NumberFormat = "@"
$ws.Cells.Item(32, 4).Value = "0.0517"
$ws.Cells.Item(32, 5).Value = "  -2.82%  "
$ws.Cells.Item(33, 5).Value = "  -4.67%  "
$ws.Cells.Item(34, 4).NumberFormat = "@"
$ws.Cells.Item(34, 4).Value = "1.82"
$ws.Cells.Item(34, 5).Value = "  -0.07%  "
$ws.Cells.Item(35, 4).Value = "1.358.57"
$ws.Cells.Item(35, 5).Value = "  -2.59%  "
$ws.Cells.Item(36, 5).Value = "  -4.63%  "
$ws.Cells.Item(37, 5).Value = "  -0.47%  "
$ws.Cells.Item(38, 4).NumberFormat = "@"
$ws.Cells.Item(38, 4).Value = "0.0188"
$ws.Cells.Item(38, 5).Value = "  -1.52%  "
$ws.Cells.Item(39, 5).Value = "  -5.43%  "
$ws.Cells.Item(40, 4).NumberFormat = "@"
$ws.Cells.Item(40, 4).Value = "2.44"
$ws.Cells.Item(40, 5).Value = "  +1.28%  "
$ws.Cells.Item(41, 2).Value = "Aave"
$ws.Cells.Item(41, 3).Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Cells.Item(41, 4).NumberFormat = "@"
$ws.Cells.Item(41, 4).Value = "81.93"
$ws.Cells.Item(41, 5).Value = "  -0.09%  "
$ws.Cells.Item(42, 2).Value = "ARBITRUM"
$ws.Cells.Item(42, 3).Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Cells.Item(42, 4).NumberFormat = "@"
$ws.Cells.Item(42, 4).Value = "0.953"
$ws.Cells.Item(42, 5).Value = "  -0.93%  "
$ws.Cells.Item(43, 5).Value = "  -1.78%  "
$ws.Cells.Item(44, 4).NumberFormat = "@"
$ws.Cells.Item(44, 4).Value = "2.80"
$ws.Cells.Item(44, 5).Value = "  -0.45%  "
$ws.Cells.Item(45, 4).NumberFormat = "@"
$ws.Cells.Item(45, 4).Value = "13.77"
$ws.Cells.Item(45, 5).Value = "  +1.53%  "
$ws.Cells.Item(46, 5).Value = "  +1.43%  "
$ws.Cells.Item(47, 4).Value = "1.973.48"
$ws.Cells.Item(47, 5).Value = "  +0.50%  "
$ws.Cells.Item(48, 4).NumberFormat = "@"
$ws.Cells.Item(48, 4).Value = "5.79"
$ws.Cells.Item(48, 5).Value = "  -4.05%  "
$ws.Cells.Item(49, 5).Value = "  -0.12%  "
$ws.Cells.Item(50, 4).NumberFormat = "@"
$ws.Cells.Item(50, 4).Value = "102.48"
$ws.Cells.Item(50, 5).Value = "  -2.37%  "
$ws.Cells.Item(51, 4).Value = "0.0₆0122"
$ws.Cells.Item(51, 5).Value = "  -4.97%  "
